$wb = $excel.ActiveWorkbook

# Update values on the "Typing" sheet
$ws = $wb.Worksheets.Item("Typing")

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 30
$ws.Range("B6").Value = 30

# Make "Typing" the active sheet/tab and set the active selection to F5
$ws.Activate()
$ws.Range("F5").Select()
